# Auto-generated edit script: updates cached numeric values in the
# "Chocobo_Profits"-style leve profit tables (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to match the refreshed market-board prices pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 2846.3103
$ws.Range("I98").Value = 1300.8235
$ws.Range("J98").Value = 5035.75
$ws.Range("K98").Value = 1300.8235
$ws.Range("L98").Value = 5035.75
$ws.Range("M98").Value = 197.1765
$ws.Range("N98").Value = -8031.75
# Row 111
$ws.Range("H111").Value = 962.58826
$ws.Range("I111").Value = 864.7778
$ws.Range("J111").Value = 1072.625
$ws.Range("K111").Value = 2594.3334
$ws.Range("L111").Value = 3217.875
$ws.Range("M111").Value = 472.6666
$ws.Range("N111").Value = -9351.875
# Row 115
$ws.Range("H115").Value = 1630
$ws.Range("I115").Value = 1577.1428
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 4731.428400000001
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -3164.428400000001
$ws.Range("N115").Value = -9134
# Row 122
$ws.Range("H122").Value = 2846.3103
$ws.Range("I122").Value = 1300.8235
$ws.Range("J122").Value = 5035.75
$ws.Range("K122").Value = 3902.4705
$ws.Range("L122").Value = 15107.25
$ws.Range("M122").Value = -1452.4705
$ws.Range("N122").Value = -20007.25
# Row 132
$ws.Range("H132").Value = 225704.28
$ws.Range("I132").Value = 3472.025
$ws.Range("J132").Value = 2003562.4
$ws.Range("K132").Value = 10416.075
$ws.Range("L132").Value = 6010687.199999999
$ws.Range("M132").Value = -7886.075000000001
$ws.Range("N132").Value = -6015747.199999999
# Row 137
$ws.Range("H137").Value = 4052.05
$ws.Range("I137").Value = 2661.1
$ws.Range("K137").Value = 7983.299999999999
$ws.Range("M137").Value = -5433.299999999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7082.408
$ws.Range("I32").Value = 6531.8438
$ws.Range("K32").Value = 6531.8438
$ws.Range("M32").Value = -6244.8438
# Row 61
$ws.Range("H61").Value = 1826.375
$ws.Range("I61").Value = 1302.2
$ws.Range("J61").Value = 2700
$ws.Range("K61").Value = 1302.2
$ws.Range("L61").Value = 2700
$ws.Range("M61").Value = -1090.2
$ws.Range("N61").Value = -3124
# Row 132
$ws.Range("H132").Value = 2128.6
$ws.Range("I132").Value = 1135.85
$ws.Range("J132").Value = 6099.6
$ws.Range("K132").Value = 3407.55
$ws.Range("L132").Value = 18298.8
$ws.Range("M132").Value = -877.5499999999997
$ws.Range("N132").Value = -23358.8
# Row 136
$ws.Range("H136").Value = 1826.375
$ws.Range("I136").Value = 1302.2
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 3906.6
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -1356.6
$ws.Range("N136").Value = -13200

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 21328.125
$ws.Range("J81").Value = 21328.125
$ws.Range("L81").Value = 21328.125
$ws.Range("N81").Value = -23450.125
# Row 84
$ws.Range("H84").Value = 21328.125
$ws.Range("J84").Value = 21328.125
$ws.Range("L84").Value = 63984.375
$ws.Range("N84").Value = -74592.375
# Row 99
$ws.Range("H99").Value = 3250.4443
$ws.Range("I99").Value = 788.625
$ws.Range("J99").Value = 5219.9
$ws.Range("K99").Value = 788.625
$ws.Range("L99").Value = 5219.9
$ws.Range("M99").Value = 709.375
$ws.Range("N99").Value = -8215.9
# Row 118
$ws.Range("H118").Value = 28890
$ws.Range("J118").Value = 28890
$ws.Range("L118").Value = 28890
$ws.Range("N118").Value = -32204
# Row 134
$ws.Range("H134").Value = 2181.3333
$ws.Range("I134").Value = 1620.138
$ws.Range("K134").Value = 4860.414
$ws.Range("M134").Value = -2325.414

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 16669670
$ws.Range("I31").Value = 1486.8667
$ws.Range("J31").Value = 33337854
$ws.Range("K31").Value = 1486.8667
$ws.Range("L31").Value = 33337854
$ws.Range("M31").Value = -1191.8667
$ws.Range("N31").Value = -33338444
# Row 34
$ws.Range("H34").Value = 16669670
$ws.Range("I34").Value = 1486.8667
$ws.Range("J34").Value = 33337854
$ws.Range("K34").Value = 1486.8667
$ws.Range("L34").Value = 33337854
$ws.Range("M34").Value = -1284.8667
$ws.Range("N34").Value = -33338258
# Row 68
$ws.Range("H68").Value = 45591.727
$ws.Range("J68").Value = 45591.727
$ws.Range("L68").Value = 45591.727
$ws.Range("N68").Value = -47089.727
# Row 71
$ws.Range("H71").Value = 45591.727
$ws.Range("J71").Value = 45591.727
$ws.Range("L71").Value = 136775.181
$ws.Range("N71").Value = -144263.181
# Row 107
$ws.Range("H107").Value = 689.9474
$ws.Range("I107").Value = 574.06665
$ws.Range("J107").Value = 1124.5
$ws.Range("K107").Value = 574.06665
$ws.Range("L107").Value = 1124.5
$ws.Range("M107").Value = 1345.93335
$ws.Range("N107").Value = -4964.5
# Row 134
$ws.Range("H134").Value = 6075.5415
$ws.Range("I134").Value = 7342.6
$ws.Range("J134").Value = 3963.7778
$ws.Range("K134").Value = 22027.8
$ws.Range("L134").Value = 11891.3334
$ws.Range("M134").Value = -19492.8
$ws.Range("N134").Value = -16961.3334

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 2202.1
$ws.Range("I70").Value = 511
$ws.Range("J70").Value = 2624.875
$ws.Range("K70").Value = 1533
$ws.Range("L70").Value = 7874.625
$ws.Range("M70").Value = -1218
$ws.Range("N70").Value = -8504.625
# Row 73
$ws.Range("H73").Value = 2202.1
$ws.Range("I73").Value = 511
$ws.Range("J73").Value = 2624.875
$ws.Range("K73").Value = 1533
$ws.Range("L73").Value = 7874.625
$ws.Range("M73").Value = -441
$ws.Range("N73").Value = -10058.625
# Row 127
$ws.Range("H127").Value = 794.6667
$ws.Range("J127").Value = 794.6667
$ws.Range("L127").Value = 2384.0001
$ws.Range("N127").Value = -12304.0001
# Row 132
$ws.Range("H132").Value = 1703.8064
$ws.Range("I132").Value = 595.26666
$ws.Range("J132").Value = 2743.0625
$ws.Range("K132").Value = 5357.39994
$ws.Range("L132").Value = 24687.5625
$ws.Range("M132").Value = -2827.39994
$ws.Range("N132").Value = -29747.5625

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 22730074
$ws.Range("I80").Value = 83335200
$ws.Range("K80").Value = 83335200
$ws.Range("M80").Value = -83334202
# Row 83
$ws.Range("H83").Value = 22730074
$ws.Range("I83").Value = 83335200
$ws.Range("K83").Value = 416676000
$ws.Range("M83").Value = -416671008

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 3913.9167
$ws.Range("I136").Value = 1595.6
$ws.Range("K136").Value = 4786.799999999999
$ws.Range("M136").Value = -2236.799999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 126325690
$ws.Range("I96").Value = 126325690
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 126325690
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -126324317
$ws.Range("N96").ClearContents()
# Row 132
$ws.Range("H132").Value = 8773495
$ws.Range("I132").Value = 982.63336
$ws.Range("K132").Value = 2947.90008
$ws.Range("M132").Value = -417.9000800000003
# Row 136
$ws.Range("H136").Value = 1676.8846
$ws.Range("I136").Value = 1281.7894
$ws.Range("J136").Value = 2749.2856
$ws.Range("K136").Value = 3845.3682
$ws.Range("L136").Value = 8247.856800000001
$ws.Range("M136").Value = -1295.3682
$ws.Range("N136").Value = -13347.8568

Write-Output "Applied scheduled Chocobo_Profits price refresh."
